$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data rows (2-78)
#    from 45192 to 45202.
$ws.Range("C2:C78").Value = 45202

# 2. Move the "A 27658-2023" record from row 7 up to row 5 (it now sorts ahead of
#    "A 62876-2018" and "A 22786-2023", which shift down to rows 6 and 7), and
#    update its data values/species text at the same time.

# Insert a blank row at position 5; this pushes the current rows 5,6,7 down to 6,7,8.
$ws.Rows.Item(5).Insert()

# The old "A 27658-2023" row is now at row 8 (its data will be rewritten into row 5
# below with the updated figures), so remove the now-duplicate row.
$ws.Rows.Item(8).Delete()

# Populate row 5 with the (updated) "A 27658-2023" record.
$ws.Range("A5").Value = "A 27658-2023"
$ws.Range("B5").Value = 45097
$ws.Range("C5").Value = 45202
$ws.Range("D5").Value = "SKÅNE LÄN"
$ws.Range("E5").Value = "ESLÖV"
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = "Havsörn`r`nSpillkråka"

$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESLOV/artfynd/A 27658-2023.xlsx", "A 27658-2023")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESLOV/kartor/A 27658-2023.png", "A 27658-2023")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESLOV/klagomål/A 27658-2023.docx", "A 27658-2023")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESLOV/klagomålsmail/A 27658-2023.docx", "A 27658-2023")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESLOV/tillsyn/A 27658-2023.docx", "A 27658-2023")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESLOV/tillsynsmail/A 27658-2023.docx", "A 27658-2023")'
